$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45271
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 45000
$ws.Range("L2").Value = 45000
$ws.Range("M2").Value = 45000
$ws.Range("P2").Value = 1800

# Row 3
$ws.Range("D3").Value = 44544
$ws.Range("H3").Value = 'Inferno'
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 35000
$ws.Range("L3").Value = 35000
$ws.Range("M3").Value = 35000
$ws.Range("P3").Value = 1400

# Row 4
$ws.Range("D4").Value = 44474
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 100000
$ws.Range("L4").Value = 100000
$ws.Range("M4").Value = 100000
$ws.Range("P4").Value = 4000

# Row 5
$ws.Range("D5").Value = 44421
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 15
$ws.Range("N5").Value = '$/caja 25 kilos'
$ws.Range("P5").Value = 3000
$ws.Range("Q5").Value = 25

# Row 6
$ws.Range("D6").Value = 44343
$ws.Range("H6").Value = 'Americana (o)'
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 36000
$ws.Range("L6").Value = 36000
$ws.Range("M6").Value = 36000
$ws.Range("P6").Value = 1440

# Row 7
$ws.Range("D7").Value = 44340
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 35000
$ws.Range("L7").Value = 35000
$ws.Range("M7").Value = 35000
$ws.Range("P7").Value = 1400

# Row 8
$ws.Range("D8").Value = 44460
$ws.Range("H8").Value = 'Americana (o)'
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 95000
$ws.Range("L8").Value = 95000
$ws.Range("M8").Value = 95000
$ws.Range("N8").Value = '$/caja 25 kilos'
$ws.Range("P8").Value = 3800
$ws.Range("Q8").Value = 25

# Row 9
$ws.Range("D9").Value = 44319
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 30000
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = 30000
$ws.Range("P9").Value = 1200

# Row 10
$ws.Range("D10").Value = 44449
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 80000
$ws.Range("L10").Value = 80000
$ws.Range("M10").Value = 80000
$ws.Range("P10").Value = 3200

# Row 11
$ws.Range("D11").Value = 44449
$ws.Range("I11").Value = 'Segunda'
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 75000
$ws.Range("L11").Value = 75000
$ws.Range("M11").Value = 75000
$ws.Range("N11").Value = '$/caja 15 kilos'
$ws.Range("P11").Value = 5000
$ws.Range("Q11").Value = 15

# Row 12
$ws.Range("D12").Value = 44553
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 45000
$ws.Range("L12").Value = 45000
$ws.Range("M12").Value = 45000
$ws.Range("P12").Value = 1800

# Row 13
$ws.Range("D13").Value = 44221
$ws.Range("J13").Value = 22
$ws.Range("K13").Value = 24000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 24545
$ws.Range("P13").Value = 982

# Row 14
$ws.Range("D14").Value = 44446
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 78000
$ws.Range("L14").Value = 78000
$ws.Range("M14").Value = 78000
$ws.Range("P14").Value = 3120

# Row 15
$ws.Range("D15").Value = 44446
$ws.Range("H15").Value = 'Inferno'
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 80000
$ws.Range("L15").Value = 80000
$ws.Range("M15").Value = 80000
$ws.Range("P15").Value = 5333

# Row 16
$ws.Range("D16").Value = 44326
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = 30000
$ws.Range("L16").Value = 30000
$ws.Range("M16").Value = 30000
$ws.Range("P16").Value = 1200

# Row 19
$ws.Range("D19").Value = 44193
$ws.Range("J19").Value = 15
$ws.Range("K19").Value = 46000
$ws.Range("L19").Value = 46000
$ws.Range("M19").Value = 46000
$ws.Range("N19").Value = '$/caja 15 kilos'
$ws.Range("P19").Value = 3067
$ws.Range("Q19").Value = 15
